$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("Month in which customer cancelled his plan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng = $d.Content
$rng.Find.Execute("Month in which customer cancelled his plan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(". If empty, it means user ")
$rng.Collapse(0)
$rng.InsertAfter("did not")
$rng.Collapse(0)
$rng.InsertAfter(" churn. Churn data is available for November, December, & January")
